$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.560.91"
$ws.Range("E2").Value = "  -7.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.174.82"
$ws.Range("E3").Value = "  -7.60%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.45"
$ws.Range("E5").Value = "  -0.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  -7.71%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "69.19"
$ws.Range("E7").Value = "  -5.57%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  -11.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.20"
$ws.Range("E10").Value = "  +3.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.32"
$ws.Range("E11").Value = "  -5.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0933"
$ws.Range("E12").Value = "  -8.69%  "

$ws.Range("E13").Value = "  -4.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.49"
$ws.Range("E14").Value = "  -9.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.493.79"
$ws.Range("E15").Value = "  -7.87%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.53"
$ws.Range("E16").Value = "  -10.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.825"
$ws.Range("E17").Value = "  -9.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.175.47"
$ws.Range("E18").Value = "  -7.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "40.461.59"
$ws.Range("E19").Value = "  -7.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0933"
$ws.Range("E20").Value = "  -9.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.96"
$ws.Range("E21").Value = "  -7.38%  "

$ws.Range("E22").Value = "  -8.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "228.20"
$ws.Range("E23").Value = "  -9.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.99"
$ws.Range("E24").Value = "  +5.98%  "

$ws.Range("E26").Value = "  -5.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.36"
$ws.Range("E27").Value = "  -5.14%  "

$ws.Range("E28").Value = "  -5.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.59"
$ws.Range("E29").Value = "  -8.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.22"
$ws.Range("E30").Value = "  -4.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.03"
$ws.Range("E31").Value = "  -10.02%  "

$ws.Range("E32").Value = "  -10.16%  "

$ws.Range("E33").Value = "  -8.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0694"
$ws.Range("E34").Value = "  -6.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.07"
$ws.Range("E35").Value = "  -5.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.54"
$ws.Range("E36").Value = "  -9.93%  "

$ws.Range("E37").Value = "  -1.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.95"
$ws.Range("E38").Value = "  +15.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.23"
$ws.Range("E39").Value = "  -7.38%  "

$ws.Range("E40").Value = "  -4.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.76"
$ws.Range("E41").Value = "  -12.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.50"
$ws.Range("E42").Value = "  -2.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.71"
$ws.Range("E43").Value = "  -14.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.57"
$ws.Range("E44").Value = "  -4.99%  "

$ws.Range("E45").Value = "  -7.11%  "

$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0973"
$ws.Range("E47").Value = "  -8.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.45"
$ws.Range("E48").Value = "  +2.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.12"
$ws.Range("E49").Value = "  +6.45%  "

$ws.Range("E50").Value = "  -7.02%  "

$ws.Range("E51").Value = "  -6.77%  "
